$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 1.003756253906252

$ws.Range("C3").Value = 1.287693099940079
$ws.Range("E3").Value = 1.224010362214401

$ws.Range("C4").Value = 0.5167526861706184
$ws.Range("E4").Value = 0.9718821796794952

$ws.Range("C5").Value = 1.187829657075357
$ws.Range("E5").Value = 1.00065194548169

$ws.Range("C6").Value = 1.071158385438342
$ws.Range("E6").Value = 0.912403143334517

$ws.Range("C7").Value = 0.9553801317191413
$ws.Range("E7").Value = 1.066801818459595

$ws.Range("C8").Value = 1.190496724073231
$ws.Range("E8").Value = 1.154811676806311

$ws.Range("C9").Value = 1.5464392869869
$ws.Range("E9").Value = 1.247870081683522

$ws.Range("C10").Value = 1.701952652941463
$ws.Range("E10").Value = 1.637918813512695

$ws.Range("C11").Value = 1.580693894992691
$ws.Range("E11").Value = 1.610567777412109

$ws.Range("C12").Value = 1.646565058924154
$ws.Range("E12").Value = 1.636439239090515

$ws.Range("C13").Value = 1.619750436871126
$ws.Range("E13").Value = 1.669486277487398

$ws.Range("C14").Value = 0.2954364073068261
$ws.Range("E14").Value = 0.8326407735962826

$ws.Range("C15").Value = -1.564297238929013
$ws.Range("E15").Value = 0.1548119563699935

$ws.Range("C16").Value = 5.937304773291885
$ws.Range("E16").Value = 2.150399152794202

$ws.Range("C17").Value = -0.2621830498131694
$ws.Range("E17").Value = 1.878976297039481

$ws.Range("C18").Value = 0.0512320434504332
$ws.Range("E18").Value = 0.2660756331863467

$ws.Range("C19").Value = 0.5998844096825495
$ws.Range("E19").Value = 0.1733734969819434
